$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New case entries for 21TRD09200 / Bunner, appended below the existing data (rows 502-525).
# Force text formatting first so numeric-looking values (e.g. "4510.111", "$ 0") are not
# auto-coerced into numbers/currency by Excel.
$ws.Range("A502:K525").NumberFormat = "@"

$newRows = @(
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','Guilty','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','Guilty','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','Guilty','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','Guilty','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','Guilty','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','Guilty','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','Guilty','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','Guilty','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','Guilty',$null,$null,$null,$null,$null),
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','No Contest','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','No Contest','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','No Contest','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','No Contest','Guilty','$ 0','$ 0',$null,$null),
    @('21TRD09200','Bunner','DUS UCM','4510.111','UCM','No Contest','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','OPERATING W/O A VALID OL - UCM','4510.12','UCM','No Contest','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS','4510.21A*','UCM','No Contest','Guilty','$ 0','$ 0','None','None'),
    @('21TRD09200','Bunner','FAILURE TO FILE REGISTRATION','4503.11','MM','No Contest','Guilty','$ 0','$ 0','None','None')
)

$startRow = 502
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}
